# Moonscraper Chart Editor Design Document - edit script
# Implements:
#  1) Split the "Mouse scroll: ..." bullet into three bullets describing the
#     new mouse controls (scroll/move chart, left click on note, left click
#     on sustain), moving the _GoBack bookmark to the end of the new last
#     bullet.
#  2) Remove the old _GoBack bookmark that used to sit after
#     "Auto-save (enable/disable)".
#  3) Shuffle the <w:lastRenderedPageBreak/> pagination marker from the
#     "Staggered" run to the "Default w/ s" run, and add one before the
#     "Screen design prototypes" picture - mirrors the repagination that
#     happens after new content is added earlier in the document.

$d = $word.ActiveDocument

$pkgNs = "xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'"
$wNs   = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Wrap-RunXml([string]$innerXml) {
    return "<pkg:package $pkgNs><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document $wNs><w:body><w:p>$innerXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# ---------------------------------------------------------------------
# 1) Rewrite the "Mouse scroll" bullet and add the two new bullets after it
# ---------------------------------------------------------------------

$mouseScrollPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Mouse scroll:*") {
        $mouseScrollPara = $p
        break
    }
}

# Replace the paragraph's text (drop the trailing paragraph mark from the range)
$bodyRange = $d.Range($mouseScrollPara.Range.Start, $mouseScrollPara.Range.End - 1)
$bodyRange.Text = "Mouse scroll: Move chart"

# Add the two new list bullets, inheriting the same list formatting
$mouseScrollPara.Range.InsertParagraphAfter()
$mouseScrollPara.Range.InsertParagraphAfter()

# Re-resolve the three paragraphs now that the document has changed shape
$mouseScrollPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Mouse scroll:*") {
        $mouseScrollPara = $p
        break
    }
}
$firstNewEnd = $mouseScrollPara.Range.End
$secondPara = $d.Range($firstNewEnd, $firstNewEnd)
$secondPara.Text = "Left click on note- move position"

$secondParaObj = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Left click on note-*") {
        $secondParaObj = $p
        break
    }
}
$thirdStart = $secondParaObj.Range.End
$thirdRange = $d.Range($thirdStart, $thirdStart)
$thirdRange.Text = "Left click on sustain- change sustain length"

# ---------------------------------------------------------------------
# 2) Move the _GoBack bookmark: delete the old one, add a new one at the
#    end of the new third bullet.
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$thirdParaObj = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Left click on sustain-*") {
        $thirdParaObj = $p
        break
    }
}
$endPos = $thirdParaObj.Range.End - 1
# A collapsed range sitting exactly at "end of paragraph text" confuses bookmark
# placement, so temporarily pad with a marker char, bookmark around it, then
# strip the marker back out - the bookmark collapses to the right spot.
$padRange = $d.Range($endPos, $endPos)
$padRange.InsertAfter("#")
$wrapRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $wrapRange)
$d.Range($endPos, $endPos + 1).Text = ""

# ---------------------------------------------------------------------
# 3) Relocate <w:lastRenderedPageBreak/> from "Staggered" to "Default w/ s",
#    and add one more before the Screen design prototype picture.
# ---------------------------------------------------------------------

# 3a. Strip it from the "Staggered" run by deleting + retyping that word
$staggeredPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Staggered*") {
        $staggeredPara = $p
        break
    }
}
$stStart = $staggeredPara.Range.Start
$stRange = $d.Range($stStart, $stStart + 9)   # "Staggered" = 9 chars
$stRange.Text = ""
$d.Range($stStart, $stStart).InsertBefore("Staggered")

# 3b. Add it before the "Default w/ s" run
$defaultPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Default w/*") {
        $defaultPara = $p
        break
    }
}
$dfStart = $defaultPara.Range.Start
$d.Range($dfStart, $dfStart).InsertXML((Wrap-RunXml("<w:r><w:lastRenderedPageBreak/></w:r>")))

# 3c. Add it before the "Screen design prototypes" picture run
$picturePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Toolpane is a moveable window*") {
        $picturePara = $p
        break
    }
}
$pStart = $picturePara.Range.Start
$d.Range($pStart, $pStart).InsertXML((Wrap-RunXml("<w:r><w:lastRenderedPageBreak/></w:r>")))
